$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new "LOAI_HD" menu entry -------------------------------------
$ws.Range("A4").Value = "LOAI_HD"
$ws.Range("D4").Value = "Loại hợp đồng"
$ws.Range("D4").Style = $ws.Range("D3").Style

# C4:C6 -> ="mnu"&A<row>  (mirrors the existing C3 pattern)
$ws.Range("C4").Formula = '="mnu"&A4'
$ws.Range("C5").Formula = '="mnu"&A5'
$ws.Range("C6").Formula = '="mnu"&A6'

# G4:G8 -> same big SQL-builder formula used by G2/G3, extended down to row 8
$gFormula = '=" IF NOT EXISTS(SELECT * FROM dbo.MENU WHERE [KEY_MENU] =  N''"&C{0}&"'' ) BEGIN INSERT INTO [dbo].[MENU]([KEY_MENU],[TEN_MENU],[TEN_MENU_ANH],[TEN_MENU_HOA],[ROOT],[HIDE],[BACK_COLOR],[IMG],[STT_MENU]) SELECT N''"&C{0}&"'' AS [KEY_MENU],    N''"&D{0}&"'' AS [TEN_MENU] ,    N''"&E{0}&"'' AS [TEN_MENU_ANH],      N''"&E{0}&"'' AS [TEN_MENU_HOA] ,[ROOT],[HIDE],[BACK_COLOR],[IMG],[STT_MENU] FROM dbo.MENU WHERE [KEY_MENU] = ''mnuDon_vi''   INSERT INTO dbo.NHOM_MENU (ID_MENU, ID_NHOM )  SELECT TOP 1 ID_MENU ,1 FROM dbo.MENU WHERE KEY_MENU =  N''"&C{0}&"'' END  "'

$ws.Range("G4").Formula = ($gFormula -f 4)
$ws.Range("G5").Formula = ($gFormula -f 5)
$ws.Range("G6").Formula = ($gFormula -f 6)
$ws.Range("G7").Formula = ($gFormula -f 7)
$ws.Range("G8").Formula = ($gFormula -f 8)

# G4:G8 need the same wrap-text cell format (style) as G1:G3
$ws.Range("G4:G8").WrapText = $true

# --- Row heights for rows 6-8 (now 90) ------------------------------------
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 90

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 18.5703125
$ws.Columns.Item(7).ColumnWidth = 50.28515625

# --- Selected cell on the sheet ---------------------------------------------
$ws.Range("E4").Select()
